$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.79659733333333
$ws.Range("H2").Value = 41.389792
$ws.Range("I2").Value = 0.9485830781324925
$ws.Range("J2").Value = 0.9485830781324925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 102.8289443333334
$ws.Range("N2").Value = 308.486833
$ws.Range("O2").Value = 0.5559120396302444
$ws.Range("P2").Value = 0.5559120396302443
$ws.Range("Q2").Value = 1418.689539178749
$ws.Range("R2").Value = 12768.20585260874
$ws.Range("S2").Value = 0.5273287537233694
$ws.Range("T2").Value = 0.5273287537233693

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.79659733333333
$ws.Range("H3").Value = 41.389792
$ws.Range("I3").Value = 0.9485830781324925
$ws.Range("J3").Value = 0.9485830781324925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3441717873742006
$ws.Range("P3").Value = 0.3441717873742006
$ws.Range("Q3").Value = 878.3276482966578
$ws.Range("R3").Value = 7904.94883466992
$ws.Range("S3").Value = 0.3264755334737809
$ws.Range("T3").Value = 0.3264755334737809

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.79659733333333
$ws.Range("H4").Value = 41.389792
$ws.Range("I4").Value = 0.9485830781324925
$ws.Range("J4").Value = 0.9485830781324925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.481835
$ws.Range("N4").Value = 55.445505
$ws.Range("O4").Value = 0.09991617299555507
$ws.Range("P4").Value = 0.09991617299555505
$ws.Range("Q4").Value = 254.9864354761067
$ws.Range("R4").Value = 2294.87791928496
$ws.Range("S4").Value = 0.09477879093534225
$ws.Range("T4").Value = 0.09477879093534224

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7478296666666666
$ws.Range("H5").Value = 2.243489
$ws.Range("I5").Value = 0.05141692186750751
$ws.Range("J5").Value = 0.05141692186750751
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 102.8289443333334
$ws.Range("N5").Value = 308.486833
$ws.Range("O5").Value = 0.5559120396302444
$ws.Range("P5").Value = 0.5559120396302443
$ws.Range("Q5").Value = 76.8985351644819
$ws.Range("R5").Value = 692.086816480337
$ws.Range("S5").Value = 0.02858328590687502
$ws.Range("T5").Value = 0.02858328590687501

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7478296666666666
$ws.Range("H6").Value = 2.243489
$ws.Range("I6").Value = 0.05141692186750751
$ws.Range("J6").Value = 0.05141692186750751
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3441717873742006
$ws.Range("P6").Value = 0.3441717873742006
$ws.Range("Q6").Value = 47.60880212564056
$ws.Range("R6").Value = 428.479219130765
$ws.Range("S6").Value = 0.01769625390041968
$ws.Range("T6").Value = 0.01769625390041968

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7478296666666666
$ws.Range("H7").Value = 2.243489
$ws.Range("I7").Value = 0.05141692186750751
$ws.Range("J7").Value = 0.05141692186750751
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.481835
$ws.Range("N7").Value = 55.445505
$ws.Range("O7").Value = 0.09991617299555507
$ws.Range("P7").Value = 0.09991617299555505
$ws.Range("Q7").Value = 13.82126450743833
$ws.Range("R7").Value = 124.391380566945
$ws.Range("S7").Value = 0.005137382060212819
$ws.Range("T7").Value = 0.005137382060212818
